$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.437127351760864
$ws.Range("B1").Value = 2.402667760848999
$ws.Range("C1").Value = 4.650745391845703
$ws.Range("D1").Value = 1.317078471183777
$ws.Range("E1").Value = 0.5113233327865601
